$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Archivo Lote")

# Remove the erroneous leading "status_pef" column (column A) that was
# mistakenly included in the batch upload template. Deleting the entire
# column shifts the remaining columns (id_peso_envio, descripcion, pmvp,
# id_servicio_franqueo) one position to the left.
$ws.Columns.Item(1).Delete()

$ws.Activate()
$ws.Columns.Item(1).Select()
